# Populate Sheet1 column A with the scenario/application comparison rows
# (each row is a single comma-joined text value, matching the source data
# that was pasted into the sheet) and leave the selection on G2 like the
# authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    "Scenario,Application,Login Method,Supported,Additional Notes",
    "1,Mobile App,Email and Password,Yes,",
    "1,Mobile App,Phone and OTP,No,",
    "1,Kiosk App,Phone and OTP,No,",
    "1,Overview,,,,Use AWS Lambda authorizer, differentiate API gateway for security, recommendation to drop MongoDB",
    "2,Mobile App,Email and Password,Yes,",
    "2,Mobile App,Phone and OTP,Yes,Future implementation required, old users can't log in, migrating changes user IDs",
    "2,Kiosk App,Phone and OTP,Yes,Mobile app users need phone verification, new users created if phone not verified",
    "2,Overview,,,,Use AWS Lambda authorizer, update client ID and secret, recommendation to drop MongoDB",
    "3,Mobile App,Email and Password,Yes,",
    "3,Mobile App,Phone and OTP,No,Requires Plus or Essentials plan",
    "3,Kiosk App,Phone and OTP,Yes,Different user pools, no user-to-transaction mapping",
    "3,Overview,,,,Use AWS Lambda authorizer, differentiate API gateway for security, recommendation to drop MongoDB"
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $rows[$i]
}

$ws.Range("G2").Select()
